$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2-305)
# from serial 45182 to serial 45184.
$ws.Range("C2:C305").Value = 45184

# Row 303 and 304 had their designation (A) and area (G) values swapped.
$a303 = $ws.Cells.Item(303, 1).Value2
$a304 = $ws.Cells.Item(304, 1).Value2
$ws.Cells.Item(303, 1).Value = $a304
$ws.Cells.Item(304, 1).Value = $a303

$g303 = $ws.Cells.Item(303, 7).Value2
$g304 = $ws.Cells.Item(304, 7).Value2
$ws.Cells.Item(303, 7).Value = $g304
$ws.Cells.Item(304, 7).Value = $g303
